$d = $word.ActiveDocument
$s = $d.Styles.Item("Normal")
try {
    Write-Output ("Font.Name=" + $s.Font.Name)
    Write-Output ("Font.Size=" + $s.Font.Size)
    Write-Output ("Font.Bold=" + $s.Font.Bold)
    Write-Output ("Font.Italic=" + $s.Font.Italic)
    Write-Output ("Font.Color=" + $s.Font.Color)
} catch {
    Write-Output "ERROR: $_"
}
try {
    Write-Output ("ParagraphFormat.LineSpacing=" + $s.ParagraphFormat.LineSpacing)
    Write-Output ("ParagraphFormat.Alignment=" + $s.ParagraphFormat.Alignment)
} catch {
    Write-Output "ERROR2: $_"
}
